# Add the "BM25_TCT" row to the "Our Results" sheet (new row 5, pushing the
# existing BM25 row down to row 6), and widen column B to fit the longer
# "Late-Int" type label.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Our Results")
$ws1.Rows.Item(5).Insert()
$ws1.Range("A5").Value = "BM25_TCT"
$ws1.Range("B5").Value = "Late-Int"
$ws1.Range("C5").Value = 0.4282
$ws1.Range("D5").Value = 0.4605
$ws1.Range("E5").Value = 0.4748
$ws1.Range("F5").Value = 0.4202
$ws1.Range("G5").Value = 0.6374
$ws1.Range("H5").Value = 0.7497
# Pick up the formatting of the (now shifted) BM25 row immediately below so
# the new row looks identical to the rest of the table.
$ws1.Range("A6:H6").Copy()
$ws1.Range("A5:H5").PasteSpecial(-4122)
$ws1.Columns.Item(2).ColumnWidth = 9.17

# Add the matching "BM25_TCT" comparison row to the "BEIR Comparison" sheet
# (new row 5, pushing the existing BM25 row down to row 6), and widen
# column B the same way.
$ws2 = $wb.Worksheets.Item("BEIR Comparison")
$ws2.Rows.Item(5).Insert()
$ws2.Range("A5").Value = "BM25_TCT"
$ws2.Range("B5").Value = "Late-Int"
$ws2.Range("C5").Value = 0.4605
$ws2.Range("D5").Value = "ColBERT"
$ws2.Range("E5").Value = 0.524
# Force this cell to stay plain text (not auto-converted to a percentage
# number) before writing the "-12.1%" label.
$ws2.Range("F5").NumberFormat = "@"
$ws2.Range("F5").Value = "-12.1%"
$ws2.Range("F5").NumberFormat = "General"
$ws2.Range("G5").Value = "Below"
# Pick up the formatting of the (now shifted) BM25 row immediately below so
# the new row looks identical to the rest of the table (red "-X%" style on
# column F, plain style elsewhere).
$ws2.Range("A6:G6").Copy()
$ws2.Range("A5:G5").PasteSpecial(-4122)
$ws2.Columns.Item(2).ColumnWidth = 9.17

$excel.CutCopyMode = 0
